# Add a new "ortho" statistics block (rows 75-80, columns D:E) to the
# hGBW worksheet, mirroring the format of the existing blocks already
# present on the sheet (e.g. the block starting at row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hGBW")

# Copy formatting from an existing header row (D3:E3) and an existing
# body row (D4:E4) onto the new block so the look matches the rest of
# the sheet.
$ws.Range("D3:E3").Copy()
$ws.Range("D75:E75").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D4:E4").Copy()
$ws.Range("D76:E80").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Row 75 - block header (graph short name)
$ws.Range("D75").Value = "Graph name"
$ws.Range("E75").Value = "ortho"

# Row 76 - description
$ws.Range("D76").Value = "Description"
$ws.Range("E76").Value = "Protein-protein orthology relations"

# Row 77
$ws.Range("D77").Value = "Number of different relations"
$ws.Range("E77").Value = "1,805,780"

# Row 79
$ws.Range("D79").Value = "Number of different organisms"
$ws.Range("E79").Value = "17"

# Row 80
$ws.Range("D80").Value = "Number of databases"
$ws.Range("E80").Value = "1 (OrthoDB)"

# Row 78 (value entered after row 80 to match the original authoring
# order reflected in the shared-strings table of the target workbook)
$ws.Range("D78").Value = "Number of different proteins"
$ws.Range("E78").Value = "139,645 (19,964 from Homo sapiens)"
